$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.727.78"
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = "'2.046.98"
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'227.42"
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('D7').Value = "'59.39"
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').Value = "'0.0834"
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = "'2.347.14"
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = "'14.39"
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = "'21.36"
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('D15').Value = "'5.47"
$ws.Range('E15').Value = '  +5.55%  '
$ws.Range('D16').Value = "'0.762"
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = "'2.037.73"
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = "'37.667.02"
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').Value = "'0.0₃0828"
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('D22').Value = "'222.37"
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  +3.52%  '
$ws.Range('D26').Value = "'169.11"
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = "'9.31"
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('E32').Value = '  +8.77%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').Value = "'0.0603"
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = "'6.46"
$ws.Range('E36').Value = '  +1.66%  '
$ws.Range('E37').Value = '  +3.87%  '
$ws.Range('E38').Value = '  +6.79%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = "'18.32"
$ws.Range('E40').Value = '  +9.47%  '
$ws.Range('D41').Value = "'1.524.88"
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').Value = "'97.90"
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = "'4.14"
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('E46').Value = '  -2.82%  '
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').Value = "'2.95"
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').Value = "'2.236.07"
$ws.Range('E51').Value = '  +0.80%  '
